$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Format column D (Price) as Text first so numeric-looking price strings
# (e.g. "0.4900", "1.000", "105.30") are preserved exactly instead of
# being auto-coerced into numbers by Excel and losing trailing zeros.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "26.462.70"
$ws.Range("E2").Value = "  -0.57%  "
$ws.Range("D3").Value = "1.715.21"
$ws.Range("E3").Value = "  -1.68%  "
$ws.Range("D4").Value = "0.9939"
$ws.Range("E4").Value = "  -0.55%  "
$ws.Range("D5").Value = "239.82"
$ws.Range("E5").Value = "  -2.93%  "
$ws.Range("D6").Value = "0.9948"
$ws.Range("E6").Value = "  -0.43%  "
$ws.Range("D7").Value = "0.4900"
$ws.Range("E7").Value = "  -0.65%  "
$ws.Range("D8").Value = "0.2590"
$ws.Range("E8").Value = "  -3.47%  "
$ws.Range("D9").Value = "0.06181"
$ws.Range("E9").Value = "  -1.75%  "
$ws.Range("D10").Value = "1.710.43"
$ws.Range("E10").Value = "  -1.85%  "
$ws.Range("D11").Value = "0.06952"
$ws.Range("E11").Value = "  -1.51%  "
$ws.Range("D12").Value = "15.57"
$ws.Range("E12").Value = "  -1.41%  "
$ws.Range("D13").Value = "0.6011"
$ws.Range("E13").Value = "  -2.65%  "
$ws.Range("D14").Value = "4.468"
$ws.Range("E14").Value = "  -2.86%  "
$ws.Range("D15").Value = "76.61"
$ws.Range("E15").Value = "  -2.20%  "
$ws.Range("D16").Value = "0.9958"
$ws.Range("E16").Value = "  -0.38%  "
$ws.Range("D17").Value = "26.295.42"
$ws.Range("E17").Value = "  -1.25%  "
$ws.Range("D18").Value = "0.9945"
$ws.Range("E18").Value = "  -0.52%  "
$ws.Range("D19").Value = "0.000007101"
$ws.Range("E19").Value = "  -2.92%  "
$ws.Range("D20").Value = "11.25"
$ws.Range("E20").Value = "  -2.91%  "
$ws.Range("D21").Value = "1.933.09"
$ws.Range("E21").Value = "  -1.44%  "
$ws.Range("D22").Value = "4.384"
$ws.Range("E22").Value = "  -4.68%  "
$ws.Range("D23").Value = "8.404"
$ws.Range("E23").Value = "  -3.93%  "
$ws.Range("D24").Value = "5.042"
$ws.Range("E24").Value = "  -4.54%  "
$ws.Range("D25").Value = "137.41"
$ws.Range("E25").Value = "  -1.64%  "
$ws.Range("D26").Value = "15.19"
$ws.Range("E26").Value = "  -1.96%  "
$ws.Range("D27").Value = "1.409"
$ws.Range("E27").Value = "  -1.01%  "
$ws.Range("D28").Value = "1.735"
$ws.Range("E28").Value = "  -2.06%  "
$ws.Range("D29").Value = "105.30"
$ws.Range("E29").Value = "  -2.34%  "
$ws.Range("D30").Value = "3.896"
$ws.Range("E30").Value = "  -4.03%  "
$ws.Range("D31").Value = "0.07963"
$ws.Range("E31").Value = "  -1.18%  "
$ws.Range("D32").Value = "3.611"
$ws.Range("E32").Value = "  -3.70%  "
$ws.Range("D33").Value = "0.04456"
$ws.Range("E33").Value = "  -3.67%  "
$ws.Range("D34").Value = "2.599"
$ws.Range("E34").Value = "  -0.41%  "
$ws.Range("D35").Value = "1.000"
$ws.Range("E35").Value = "  -2.04%  "
$ws.Range("D36").Value = "0.6157"
$ws.Range("E36").Value = "  -3.84%  "
$ws.Range("D37").Value = "0.9562"
$ws.Range("E37").Value = "  +6.09%  "
$ws.Range("D38").Value = "1.996"
$ws.Range("E38").Value = "  -4.52%  "
$ws.Range("D39").Value = "2.362"
$ws.Range("E39").Value = "  -2.68%  "
$ws.Range("D40").Value = "0.9947"
$ws.Range("E40").Value = "  -0.84%  "
$ws.Range("D41").Value = "0.01478"
$ws.Range("E41").Value = "  -1.89%  "
$ws.Range("D42").Value = "99.61"
$ws.Range("E42").Value = "  -2.16%  "
$ws.Range("D43").Value = "5.409"
$ws.Range("E43").Value = "  -0.41%  "
$ws.Range("D44").Value = "0.3811"
$ws.Range("E44").Value = "  -3.21%  "
$ws.Range("D45").Value = "6.854"
$ws.Range("E45").Value = "  -0.92%  "
$ws.Range("D46").Value = "0.1150"
$ws.Range("E46").Value = "  -3.25%  "
$ws.Range("D47").Value = "0.05341"
$ws.Range("E47").Value = "  -1.08%  "
$ws.Range("D48").Value = "30.36"
$ws.Range("E48").Value = "  -0.84%  "
$ws.Range("D49").Value = "7.688"
$ws.Range("E49").Value = "  -2.25%  "
$ws.Range("D50").Value = "51.10"
$ws.Range("E50").Value = "  -1.44%  "
$ws.Range("D51").Value = "0.9979"
$ws.Range("E51").Value = "  -0.45%  "

# Restore the default style on column D so formatting matches the original
# workbook (only the text values should differ).
$ws.Range("D2:D51").Style = "Normal"
